# Fixed bug in array expansion
#
# Populate the lookup table used by an INDEX/MATCH array formula on the
# "Referencing" sheet, and enter that array formula itself, spilling its
# expanded results down column D.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Referencing")

# Header row for the lookup table (row 15, columns C:F)
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 4

# Data rows for the lookup table (rows 16-19, columns C:F)
$ws.Range("C16").Value = 1.4535833325868115
$ws.Range("D16").Value = 1.4535833325868115
$ws.Range("E16").Value = 1.5117266658902839
$ws.Range("F16").Value = 1.5407983325420203

$ws.Range("C17").Value = 9.0545454545454547
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = 18
$ws.Range("F17").Value = 18

$ws.Range("C18").Value = 0.36811506356713858
$ws.Range("D18").Value = 0.36811506356713858
$ws.Range("E18").Value = 0.40588480110308967
$ws.Range("F18").Value = 0.42190146532760275

$ws.Range("C19").Value = 0.65100000000000002
$ws.Range("D19").Value = 0.65100000000000002
$ws.Range("E19").Value = 0.65100000000000002
$ws.Range("F19").Value = 0.65100000000000002

# Lookup key and the (fixed) array formula that expands down D22:D25
$ws.Range("C22").Value = 4
$ws.Range("D22:D25").FormulaArray = "=INDEX(C16:F19,,MATCH(C22,C15:F15,0))"

# Leave the selection where the author left it after fixing the bug
$ws.Range("C23").Select()
